$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 644, pushing existing rows 644-692 down to 647-695.
$ws.Range("644:646").Insert()

# Shared template values for these Brocoli / Agricola del Norte S.A. de Arica rows.
$mercadoId = 1
$mercado = "Agrícola del Norte S.A. de Arica"
$region = "Arica y Parinacota"
$codreg = 15
$categoriaId = 100112023
$categoria = "Brócoli"
$variedad = "Sin especificar"
$unidad = "`$/unidad"
$origen = "Región de Arica y Parinacota"
$kgUnidades = 1
$clasificacion = "Hortaliza"

function Set-BrocoliRow($RowNum, $Fecha, $Calidad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm) {
    $ws.Cells.Item($RowNum, 1).Value = $mercadoId
    $ws.Cells.Item($RowNum, 2).Value = $mercado
    $ws.Cells.Item($RowNum, 3).Value = $region
    $ws.Cells.Item($RowNum, 4).Value = $Fecha
    $ws.Cells.Item($RowNum, 5).Value = $codreg
    $ws.Cells.Item($RowNum, 6).Value = $categoriaId
    $ws.Cells.Item($RowNum, 7).Value = $categoria
    $ws.Cells.Item($RowNum, 8).Value = $variedad
    $ws.Cells.Item($RowNum, 9).Value = $Calidad
    $ws.Cells.Item($RowNum, 10).Value = $Volumen
    $ws.Cells.Item($RowNum, 11).Value = $PrecioMin
    $ws.Cells.Item($RowNum, 12).Value = $PrecioMax
    $ws.Cells.Item($RowNum, 13).Value = $PrecioProm
    $ws.Cells.Item($RowNum, 14).Value = $unidad
    $ws.Cells.Item($RowNum, 15).Value = $origen
    $ws.Cells.Item($RowNum, 16).Value = $PrecioProm
    $ws.Cells.Item($RowNum, 17).Value = $kgUnidades
    $ws.Cells.Item($RowNum, 18).Value = $clasificacion
}

Set-BrocoliRow 644 45265 "Primera" 1200 600 700 650
Set-BrocoliRow 645 45265 "Segunda" 1200 450 500 475
Set-BrocoliRow 646 45265 "Tercera" 1200 350 400 375
